$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "rotate platform (90 mode)" event row (row 13) ---
# Event type id
$ws.Range("A13").Value2 = 14
# Event description
$ws.Range("B13").Value2 = "rotate platform (90 mode)"
# Value description
$ws.Range("C13").Value2 = "see rotation table below"

# Match formatting of the existing table rows:
#  - column A uses the grey/right-aligned numeric style (same as A3:A12)
#  - column C uses the italic "see ... table below" style (same as C12)
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update the view state (scroll position / selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("F12").Select()
